$d = $word.ActiveDocument

# --- 1. Locate the paragraph that currently holds only the page-break run.
#     It sits right after the (mostly empty) "Alex: Describe the system"
#     answer paragraph, and right before the next section's heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Length -ge 1 -and [int][char]$t[0] -eq 12) {
        $target = $p
        break
    }
}

$r = $target.Range
$r.Collapse(1)

$paraA = "The software will be a restaurant management. It should allow for the tracking and processing of multiple restaurant functions. The software should allow an employee to sign in to the terminal and clock in for the day. Upon clocking in the employee should be assigned a section of tables to work. While the employee is signed in to the terminal they should be able to place orders for each table and order any existing orders. The software should also allow for the tracking of inventory numbers so that when an employee places an order the inventory is updated automatically to represent the usage of materials. The software should also allow for the generation of a bill for each table when requested. Finally the employee should be able to clock out at the end of the night." + " The employee should be able to sign out of the terminal at any time allowing someone else to sign in while maintaining their table information." + " Upon the closing of the software a report should be generated for the day that includes sales totals as well as inventory use."

$paraB = "All of the software should be used through a graphic user interface that allows for using a mouse for input. The menus displayed should appropriately separate the functions by category such as sign in/sign out, place order, generate bill and log out."

# --- 2. Type the first new paragraph's text (tab + long description).
$r.InsertBefore("`t" + $paraA)

# --- 3. Split here: the page-break run stays behind, in a brand-new
#     paragraph, right after this point.
$r.Collapse(0)
$r.InsertParagraphAfter()

# --- 4. Re-fetch the (now second) paragraph that holds the page-break run
#     and type the second new paragraph's text before it.
$idx = $target.Range.Information(1)
$p2 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Length -ge 1 -and [int][char]$t[0] -eq 12) {
        $p2 = $p
        break
    }
}
$r2 = $p2.Range
$r2.Collapse(1)
$r2.InsertBefore("`t" + $paraB)

# --- 5. Move the _GoBack bookmark: remove it from its old location (a table
#     cell) and re-create it right before the page-break run, i.e. right
#     after the text we just typed in the second new paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)
